$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Fill in the "Front End Design" task for the three existing rows (48-50
#    worth of attendance, sheet rows 50-52) that previously had no task text
#    in column D.
# ---------------------------------------------------------------------------
$ws.Range("D50").Value = "Front End Design"
$ws.Range("D51").Value = "Front End Design"
$ws.Range("D52").Value = "Front End Design"

# ---------------------------------------------------------------------------
# 2) Touch each brand-new task string once, in the same order the source
#    workbook first introduces them, so the shared-string table comes out
#    in the expected order (Front End Design, Validation, Dropdown
#    onloading, "Dropdown onloading " w/ trailing space, Client screen
#    connection). The real values land permanently a few lines down when
#    the full rows are populated.
# ---------------------------------------------------------------------------
$ws.Range("D61").Value = "Validation"
$ws.Range("D56").Value = "Dropdown onloading"
$ws.Range("D60").Value = "Dropdown onloading "
$ws.Range("D71").Value = "Client screen connection"

# ---------------------------------------------------------------------------
# 3) Copy the number format already used for the date column (B48 carries
#    style s="1" -> built-in date numFmtId 14) down onto the new date cells
#    B53:B71 so they reuse the existing style instead of Excel minting a
#    new custom number format entry.
# ---------------------------------------------------------------------------
$ws.Range("B48").Copy()
$ws.Range("B53:B71").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 4) Populate the 19 new attendance rows (53-71) with their final values.
# ---------------------------------------------------------------------------
$ws.Range("A53").Value = 51
$ws.Range("B53").Value = 43133
$ws.Range("C53").Value = "Saturday"
$ws.Range("D53").Value = "PLACEMENT DRIVE"
$ws.Range("A54").Value = 52
$ws.Range("B54").Value = 43134
$ws.Range("C54").Value = "Sunday"
$ws.Range("D54").Value = "HOLIDAY"
$ws.Range("A55").Value = 53
$ws.Range("B55").Value = 43135
$ws.Range("C55").Value = "Sunday"
$ws.Range("D55").Value = "HOLIDAY"
$ws.Range("A56").Value = 54
$ws.Range("B56").Value = 43136
$ws.Range("C56").Value = "Monday"
$ws.Range("D56").Value = "Dropdown onloading"
$ws.Range("A57").Value = 55
$ws.Range("B57").Value = 43137
$ws.Range("C57").Value = "Tuesday"
$ws.Range("D57").Value = "Dropdown onloading"
$ws.Range("A58").Value = 56
$ws.Range("B58").Value = 43138
$ws.Range("C58").Value = "Wednesday"
$ws.Range("D58").Value = "Dropdown onloading"
$ws.Range("A59").Value = 57
$ws.Range("B59").Value = 43139
$ws.Range("C59").Value = "Thursday"
$ws.Range("D59").Value = "Dropdown onloading"
$ws.Range("A60").Value = 58
$ws.Range("B60").Value = 43140
$ws.Range("C60").Value = "Friday"
$ws.Range("D60").Value = "Dropdown onloading "
$ws.Range("A61").Value = 59
$ws.Range("B61").Value = 43141
$ws.Range("C61").Value = "Saturday"
$ws.Range("D61").Value = "Validation"
$ws.Range("A62").Value = 60
$ws.Range("B62").Value = 43142
$ws.Range("C62").Value = "Sunday"
$ws.Range("D62").Value = "HOLIDAY"
$ws.Range("A63").Value = 61
$ws.Range("B63").Value = 43143
$ws.Range("C63").Value = "Monday"
$ws.Range("D63").Value = "Validation"
$ws.Range("A64").Value = 62
$ws.Range("B64").Value = 43144
$ws.Range("C64").Value = "Tuesday"
$ws.Range("D64").Value = "Validation"
$ws.Range("A65").Value = 63
$ws.Range("B65").Value = 43145
$ws.Range("C65").Value = "Wednesday"
$ws.Range("D65").Value = "Validation"
$ws.Range("A66").Value = 64
$ws.Range("B66").Value = 43146
$ws.Range("C66").Value = "Thursday"
$ws.Range("D66").Value = "Validation"
$ws.Range("A67").Value = 65
$ws.Range("B67").Value = 43147
$ws.Range("C67").Value = "Friday"
$ws.Range("D67").Value = "PLACEMENT DRIVE"
$ws.Range("A68").Value = 66
$ws.Range("B68").Value = 43148
$ws.Range("C68").Value = "Saturday"
$ws.Range("D68").Value = "PLACEMENT DRIVE"
$ws.Range("A69").Value = 67
$ws.Range("B69").Value = 43149
$ws.Range("C69").Value = "Sunday"
$ws.Range("D69").Value = "HOLIDAY"
$ws.Range("A70").Value = 68
$ws.Range("B70").Value = 43150
$ws.Range("C70").Value = "Monday"
$ws.Range("D70").Value = "PLACEMENT DRIVE"
$ws.Range("A71").Value = 69
$ws.Range("B71").Value = 43151
$ws.Range("C71").Value = "Tuesday"
$ws.Range("D71").Value = "Client screen connection"

# ---------------------------------------------------------------------------
# 5) Move the selection to the last entered cell, matching the author's
#    final cursor position, and scroll the view back to the top of the
#    sheet (the saved view no longer has a frozen topLeftCell below A1).
# ---------------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null
$ws.Range("D71").Select() | Out-Null
